$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-16 Sunday", "2025-02-17 Monday"),
    @("517×3=", "955×2="),
    @("626×8=", "977×4="),
    @("648×8=", "820×7="),
    @("888×9=", "553×9="),
    @("197×7=", "827×4="),
    @("333×6=", "482×7="),
    @("803×6=", "470×4="),
    @("722×5=", "375×7="),
    @("779×3=", "591×9="),
    @("498×5=", "430×3="),
    @("161×5=", "842×2="),
    @("680×8=", "975×7="),
    @("612×8=", "662×3="),
    @("237×5=", "579×5="),
    @("356×6=", "438×4="),
    @("252×7=", "914×3="),
    @("164×4=", "353×7="),
    @("503×4=", "456×6="),
    @("353×9=", "626×6="),
    @("884×3=", "250×8="),
    @("233×6=", "927×4="),
    @("728×3=", "266×3="),
    @("125×3=", "492×3="),
    @("782×4=", "877×5="),
    @("631×3=", "457×7=")
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

$d.Save()
